$d = $word.ActiveDocument

# --- 1. Create the three new character styles ------------------------------

$GaNStyle = $d.Styles.Add("GaNStyle", 2)
$GaNStyle.Font.Name = "Calibri"
$GaNStyle.Font.NameAscii = "Calibri"
$GaNStyle.Font.Size = 14

$GaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$GaNParagraph.Font.Name = "Calibri"
$GaNParagraph.Font.NameAscii = "Calibri"
$GaNParagraph.Font.Size = 10

$GaNLinks = $d.Styles.Add("GaNLinks", 2)
$GaNLinks.Font.Name = "Calibri"
$GaNLinks.Font.NameAscii = "Calibri"
$GaNLinks.Font.Bold = $true
$GaNLinks.Font.Color = 8388608
$GaNLinks.Font.Size = 9.5
$GaNLinks.Font.Underline = 1

# --- 2. The four "Pegasuksen tähdistö ..." runs: add a trailing period and
#        apply the GaNStyle character style -------------------------------

$rng = $d.Content
$rng.Find.ClearFormatting()
while ($rng.Find.Execute("Pegasuksen tähdistö havainnointijaksot vuonna 2022: 8.-17.10., 7.-16.11.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Text = "Pegasuksen tähdistö havainnointijaksot vuonna 2022: 8.-17.10., 7.-16.11.."
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- 3. The "Osallistut maailmanlaajuiseen ..." run: apply GaNParagraph ----

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
if ($rng2.Find.Execute("Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Pegasuksen tähdistö miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
}

# --- 4. The "Tämän oppaan kartat piirsi ..." run: apply GaNLinks -----------

$rng3 = $d.Content
$rng3.Find.ClearFormatting()
if ($rng3.Find.Execute("Tämän oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
}
